$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.619.10"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "1.635.84"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.42"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("E6").Value = "  -1.36%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.92"
$ws.Range("E8").Value = "  -0.74%  "
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").Value = "1.868.00"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").Value = "1.628.10"
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("E14").Value = "  -0.13%  "
$ws.Range("E15").Value = "  -1.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.51"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").Value = "27.612.82"
$ws.Range("E17").Value = "  +0.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.21"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.74"
$ws.Range("E19").Value = "  +1.90%  "
$ws.Range("D20").Value = "0.0₃0721"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.01"
$ws.Range("E23").Value = "  +4.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.26"
$ws.Range("E25").Value = "  +2.12%  "
$ws.Range("E26").Value = "  -1.15%  "
$ws.Range("E27").Value = "  -1.48%  "
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").Value = "1.453.98"
$ws.Range("E33").Value = "  +2.77%  "
$ws.Range("E34").Value = "  -1.55%  "
$ws.Range("E35").Value = "  -0.75%  "
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.875"
$ws.Range("E38").Value = "  -1.07%  "
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.892"
$ws.Range("E40").Value = "  +8.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.70"
$ws.Range("E41").Value = "  +7.96%  "
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("E43").Value = "  -0.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.59"
$ws.Range("E44").Value = "  +1.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.47"
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").Value = "1.777.92"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.70"
$ws.Range("E48").Value = "  +1.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.19"
$ws.Range("E49").Value = "  -2.10%  "
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.78"
$ws.Range("E51").Value = "  +0.53%  "
